$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3501.3333
$ws.Range("I51").Value = 3500
$ws.Range("K51").Value = 3500
$ws.Range("M51").Value = -3016

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 266.9
$ws.Range("I53").Value = 126.5
$ws.Range("J53").Value = 594.5
$ws.Range("K53").Value = 126.5
$ws.Range("L53").Value = 594.5
$ws.Range("M53").Value = 510.5
$ws.Range("N53").Value = -1868.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4043.2632
$ws.Range("I86").Value = 1967.6666
$ws.Range("J86").Value = 4432.4375
$ws.Range("K86").Value = 1967.6666
$ws.Range("L86").Value = 4432.4375
$ws.Range("M86").Value = -844.6666
$ws.Range("N86").Value = -6678.4375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4043.2632
$ws.Range("I89").Value = 1967.6666
$ws.Range("J89").Value = 4432.4375
$ws.Range("K89").Value = 9838.333000000001
$ws.Range("L89").Value = 22162.1875
$ws.Range("M89").Value = -4222.333000000001
$ws.Range("N89").Value = -33394.1875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 62499.2
$ws.Range("I95").Value = 30000
$ws.Range("K95").Value = 30000
$ws.Range("M95").Value = -27254

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2141.3572
$ws.Range("I106").Value = 1559.875
$ws.Range("K106").Value = 1559.875
$ws.Range("M106").Value = -928.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3037.4517
$ws.Range("I138").Value = 2336.5293
$ws.Range("J138").Value = 3888.5715
$ws.Range("K138").Value = 7009.5879
$ws.Range("L138").Value = 11665.7145
$ws.Range("M138").Value = -1869.5879
$ws.Range("N138").Value = -21945.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 264998.75
$ws.Range("I6").Value = 19995
$ws.Range("K6").Value = 19995
$ws.Range("M6").Value = -19822

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 4752508.5
$ws.Range("I8").Value = 14250000
$ws.Range("J8").Value = 3762.5
$ws.Range("K8").Value = 14250000
$ws.Range("L8").Value = 3762.5
$ws.Range("M8").Value = -14249856
$ws.Range("N8").Value = -4050.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 668566.7
$ws.Range("J11").Value = 2850
$ws.Range("L11").Value = 2850
$ws.Range("N11").Value = -3138

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3793080.8
$ws.Range("I32").Value = 4697595
$ws.Range("J32").Value = 15404.647
$ws.Range("K32").Value = 4697595
$ws.Range("L32").Value = 15404.647
$ws.Range("M32").Value = -4697308
$ws.Range("N32").Value = -15978.647

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4330.143
$ws.Range("I45").Value = 4740.5713
$ws.Range("J45").Value = 3509.2856
$ws.Range("K45").Value = 4740.5713
$ws.Range("L45").Value = 3509.2856
$ws.Range("M45").Value = -4363.5713
$ws.Range("N45").Value = -4263.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1863445.5
$ws.Range("I61").Value = 1863445.5
$ws.Range("K61").Value = 1863445.5
$ws.Range("M61").Value = -1863233.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 56799.883
$ws.Range("J62").Value = 56799.883
$ws.Range("L62").Value = 56799.883
$ws.Range("N62").Value = -58047.883

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5776.8887
$ws.Range("J63").Value = 7121.923
$ws.Range("L63").Value = 7121.923
$ws.Range("N63").Value = -8493.922999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 56799.883
$ws.Range("J65").Value = 56799.883
$ws.Range("L65").Value = 170399.649
$ws.Range("N65").Value = -176639.649

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5776.8887
$ws.Range("J66").Value = 7121.923
$ws.Range("L66").Value = 35609.615
$ws.Range("N66").Value = -42473.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2110.2
$ws.Range("I88").Value = 1950
$ws.Range("J88").Value = 2350.5
$ws.Range("K88").Value = 1950
$ws.Range("L88").Value = 2350.5
$ws.Range("M88").Value = -1544
$ws.Range("N88").Value = -3162.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2110.2
$ws.Range("I91").Value = 1950
$ws.Range("J91").Value = 2350.5
$ws.Range("K91").Value = 1950
$ws.Range("L91").Value = 2350.5
$ws.Range("M91").Value = -546
$ws.Range("N91").Value = -5158.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 33805.168
$ws.Range("J96").Value = 33805.168
$ws.Range("L96").Value = 33805.168
$ws.Range("N96").Value = -39297.168

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5465.3335
$ws.Range("I102").Value = 4456.857
$ws.Range("K102").Value = 4456.857
$ws.Range("M102").Value = -2834.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 557215.2
$ws.Range("I132").Value = 648512.3
$ws.Range("K132").Value = 1945536.9
$ws.Range("M132").Value = -1943006.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1863445.5
$ws.Range("I136").Value = 1863445.5
$ws.Range("K136").Value = 5590336.5
$ws.Range("M136").Value = -5587786.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 866.6667
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 866.6667
$ws.Range("K11").Value = 0
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = 866.6667
$ws.Range("N11").Value = -1146.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000.25
$ws.Range("I86").Value = 1667
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1667
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -544
$ws.Range("N86").Value = -5246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2000.25
$ws.Range("I89").Value = 1667
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 8335
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -2719
$ws.Range("N89").Value = -26232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2038
$ws.Range("I105").Value = 1992
$ws.Range("K105").Value = 1992
$ws.Range("M105").Value = -245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2028.7778
$ws.Range("I6").Value = 1969.875
$ws.Range("K6").Value = 1969.875
$ws.Range("M6").Value = -1856.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 100208.8
$ws.Range("J119").Value = 100208.8
$ws.Range("L119").Value = 100208.8
$ws.Range("N119").Value = -109884.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5495
$ws.Range("I80").Value = 4995
$ws.Range("J80").Value = 5995
$ws.Range("K80").Value = 14985
$ws.Range("L80").Value = 17985
$ws.Range("M80").Value = -14049
$ws.Range("N80").Value = -19857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5495
$ws.Range("I83").Value = 4995
$ws.Range("J83").Value = 5995
$ws.Range("K83").Value = 44955
$ws.Range("L83").Value = 53955
$ws.Range("M83").Value = -40275
$ws.Range("N83").Value = -63315

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1747.9333
$ws.Range("J113").Value = 1588.2727
$ws.Range("L113").Value = 4764.8181
$ws.Range("N113").Value = -9104.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2591.6
$ws.Range("I132").Value = 1671.1428
$ws.Range("J132").Value = 3087.2307
$ws.Range("K132").Value = 15040.2852
$ws.Range("L132").Value = 27785.0763
$ws.Range("M132").Value = -12510.2852
$ws.Range("N132").Value = -32845.0763

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 236.3
$ws.Range("I2").Value = 130.28572
$ws.Range("K2").Value = 130.28572
$ws.Range("M2").Value = -17.28572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2316000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 20676.291
$ws.Range("I24").Value = 10920.6
$ws.Range("J24").Value = 27644.643
$ws.Range("K24").Value = 10920.6
$ws.Range("L24").Value = 27644.643
$ws.Range("M24").Value = -10747.6
$ws.Range("N24").Value = -27990.643

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8290.714
$ws.Range("I70").Value = 9349.200000000001
$ws.Range("J70").Value = 5644.5
$ws.Range("K70").Value = 9349.200000000001
$ws.Range("L70").Value = 5644.5
$ws.Range("M70").Value = -9079.200000000001
$ws.Range("N70").Value = -6184.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8290.714
$ws.Range("I73").Value = 9349.200000000001
$ws.Range("J73").Value = 5644.5
$ws.Range("K73").Value = 9349.200000000001
$ws.Range("L73").Value = 5644.5
$ws.Range("M73").Value = -8413.200000000001
$ws.Range("N73").Value = -7516.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2763.3333
$ws.Range("I97").Value = 341.5
$ws.Range("J97").Value = 7607
$ws.Range("K97").Value = 341.5
$ws.Range("L97").Value = 7607
$ws.Range("M97").Value = 154.5
$ws.Range("N97").Value = -8599

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10199.667
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 14799.5
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 14799.5
$ws.Range("M99").Value = 1246
$ws.Range("N99").Value = -19291.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2489.4285
$ws.Range("I132").Value = 1489
$ws.Range("J132").Value = 4990.5
$ws.Range("K132").Value = 4467
$ws.Range("L132").Value = 14971.5
$ws.Range("M132").Value = -1937
$ws.Range("N132").Value = -20031.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3253.92
$ws.Range("I46").Value = 1416.6666
$ws.Range("J46").Value = 3834.1052
$ws.Range("K46").Value = 1416.6666
$ws.Range("L46").Value = 3834.1052
$ws.Range("M46").Value = -1228.6666
$ws.Range("N46").Value = -4210.1052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3742.6428
$ws.Range("I61").Value = 2399.875
$ws.Range("K61").Value = 2399.875
$ws.Range("M61").Value = -2197.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3742.6428
$ws.Range("I113").Value = 2399.875
$ws.Range("K113").Value = 2399.875
$ws.Range("M113").Value = -229.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4533.6875
$ws.Range("I122").Value = 4400.0415
$ws.Range("J122").Value = 4934.625
$ws.Range("K122").Value = 13200.1245
$ws.Range("L122").Value = 14803.875
$ws.Range("M122").Value = -10750.1245
$ws.Range("N122").Value = -19703.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5213.0454
$ws.Range("I136").Value = 4889.857
$ws.Range("K136").Value = 14669.571
$ws.Range("M136").Value = -12119.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 6070
$ws.Range("I51").Value = 6070
$ws.Range("K51").Value = 6070
$ws.Range("M51").Value = -5560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5250.25
$ws.Range("I100").Value = 5999.6665
$ws.Range("K100").Value = 11999.333
$ws.Range("M100").Value = -11458.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 45510.43
$ws.Range("J125").Value = 45510.43
$ws.Range("L125").Value = 45510.43
$ws.Range("N125").Value = -55350.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4068870.8
$ws.Range("I132").Value = 5379503
$ws.Range("K132").Value = 16138509
$ws.Range("M132").Value = -16135979
